$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-44 down to 4-45
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new data point
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44817
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100112052
$ws.Cells.Item(3, 7).Value = "Albahaca"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 1300
$ws.Cells.Item(3, 12).Value = 1500
$ws.Cells.Item(3, 13).Value = 1400
$ws.Cells.Item(3, 14).Value = "`$/paquete"
$ws.Cells.Item(3, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 16).Value = 1400
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same number format/style as the other date cells in column D
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
